# Scene 82 touch-up:
#  1. Merge the three split runs that make up the "Mom (...)..." sincere
#     line back into a single run (no visible text change, just undoing
#     the run fragmentation).
#  2. Remove the stray "Mom: (neutral raised_eyebrow):" placeholder
#     paragraph entirely.

$d = $word.ActiveDocument

# --- 1. Merge the fragmented "Mom (neutral sincere): I'm really sorry..." run ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*neutral sincere*") {
        $start = $p.Range.Start
        $end = $p.Range.End
        $full = $d.Range($start, $end)
        $full.Text = "Mom (neutral sincere): I’m really sorry. When you left this morning, I assumed you’d be going to school."
        break
    }
}

# --- 2. Delete the "Mom: (neutral raised_eyebrow):" paragraph entirely ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Mom: (neutral raised_eyebrow):*") {
        $p.Range.Delete()
        break
    }
}
